$d = $word.ActiveDocument
$sec = $d.Sections(1)
$sec.PageSetup.OddAndEvenPagesHeaderFooter = 1
$sec.PageSetup.DifferentFirstPageHeaderFooter = 1

$h1 = $sec.Headers(1); $h1.Range.Text = ""
$h2 = $sec.Headers(2); $h2.Range.Text = ""
$h3 = $sec.Headers(3); $h3.Range.Text = ""
$f1 = $sec.Footers(1)
$f2 = $sec.Footers(2); $f2.Range.Text = ""
$f3 = $sec.Footers(3); $f3.Range.Text = ""

$cr = [char]13
$f1.Range.Text = "X" + $cr

$p1 = $f1.Range.Paragraphs(1).Range
$fld = $d.Fields.Add($p1, 33)
Write-Host ("after add: [" + $f1.Range.Text + "]")

$found = $f1.Range.Find.Execute("X", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
Write-Host ("Find replaced: " + $found)
Write-Host ("final: [" + $f1.Range.Text + "]")
Write-Host ("Paragraphs.Count=" + $f1.Range.Paragraphs.Count)
